# test solution to #226
# - Switch the active sheet from "drop" (3rd tab) to "numeric" (2nd tab)
# - Update the selection on the "numeric" sheet to E3
# - Fix the value in E2 on the "numeric" sheet from " " to "*"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")

# Correct the flagged value
$ws.Range("E2").Value = "*"

# Make "numeric" the active sheet and move the selection to E3
$ws.Activate()
$ws.Range("E3").Select()
